{"js": "// Fix typo in question 1b: \"FROM. \" -> \"FROM \" (stray period after FROM).\n// The paragraph runs \"...FRO\" + \"M\" + \".\" + \" \" get collapsed so the\n// period is dropped while the bold \"M \" text is preserved as one run.\nconst body = context.document.body;\n\nconst results = body.search(\"M. \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the typo text 'M. ' to fix.\");\n}\n\n// Replace the matched range's text, dropping the stray period.\nresults.items[0].insertText(\"M \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix typo in question 1b: remove the stray period in \"...FROM. \" so it\n# reads \"...FROM \" instead.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"M. \"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    # $range now spans exactly \"M. \" (the bold \"M\", the stray period, and\n    # the trailing space). The period is the second character - delete it\n    # in place so \"M\" and \" \" stay put and simply merge into one run.\n    $periodRange = $d.Range($range.Start + 1, $range.Start + 2)\n    if ($periodRange.Text -eq \".\") {\n        $periodRange.Delete()\n    }\n}\n"}
